$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 50.70817566666667
$ws.Range("H2").Value = 152.124527
$ws.Range("I2").Value = 0.5661129211027078
$ws.Range("J2").Value = 0.5661129211027077
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.778439
$ws.Range("N2").Value = 11.335317
$ws.Range("O2").Value = 0.4252971528324392
$ws.Range("P2").Value = 0.4252971528324392
$ws.Range("Q2").Value = 191.5977485577843
$ws.Range("R2").Value = 1724.379737020059
$ws.Range("S2").Value = 0.240766213526637
$ws.Range("T2").Value = 0.2407662135266369
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 50.70817566666667
$ws.Range("H3").Value = 152.124527
$ws.Range("I3").Value = 0.5661129211027078
$ws.Range("J3").Value = 0.5661129211027077
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("N3").Value = 13.00021
$ws.Range("O3").Value = 0.4877633593505858
$ws.Range("P3").Value = 0.4877633593505858
$ws.Range("Q3").Value = 219.7389774611856
$ws.Range("R3").Value = 1977.65079715067
$ws.Range("S3").Value = 0.2761291401688299
$ws.Range("T3").Value = 0.2761291401688299
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 50.70817566666667
$ws.Range("H4").Value = 152.124527
$ws.Range("I4").Value = 0.5661129211027078
$ws.Range("J4").Value = 0.5661129211027077
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2909853333333334
$ws.Range("N4").Value = 0.8729560000000001
$ws.Range("O4").Value = 0.03275300561492853
$ws.Range("P4").Value = 0.03275300561492853
$ws.Range("Q4").Value = 14.75533539909022
$ws.Range("R4").Value = 132.798018591812
$ws.Range("S4").Value = 0.01854189968356058
$ws.Range("T4").Value = 0.01854189968356058
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 50.70817566666667
$ws.Range("H5").Value = 152.124527
$ws.Range("I5").Value = 0.5661129211027078
$ws.Range("J5").Value = 0.5661129211027077
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4814053333333333
$ws.Range("N5").Value = 1.444216
$ws.Range("O5").Value = 0.0541864822020464
$ws.Range("P5").Value = 0.05418648220204641
$ws.Range("Q5").Value = 24.41118620953689
$ws.Range("R5").Value = 219.700675885832
$ws.Range("S5").Value = 0.03067566772368038
$ws.Range("T5").Value = 0.03067566772368038
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.08683666666667
$ws.Range("H6").Value = 51.26051
$ws.Range("I6").Value = 0.1907597520636141
$ws.Range("J6").Value = 0.1907597520636141
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.778439
$ws.Range("N6").Value = 11.335317
$ws.Range("O6").Value = 0.4252971528324392
$ws.Range("P6").Value = 0.4252971528324392
$ws.Range("Q6").Value = 64.56157004796333
$ws.Range("R6").Value = 581.05413043167
$ws.Range("S6").Value = 0.08112957942767708
$ws.Range("T6").Value = 0.08112957942767708
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.08683666666667
$ws.Range("H7").Value = 51.26051
$ws.Range("I7").Value = 0.1907597520636141
$ws.Range("J7").Value = 0.1907597520636141
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("N7").Value = 13.00021
$ws.Range("O7").Value = 0.4877633593505858
$ws.Range("P7").Value = 0.4877633593505858
$ws.Range("Q7").Value = 74.04415496745555
$ws.Range("R7").Value = 666.3973947071
$ws.Range("S7").Value = 0.09304561749543323
$ws.Range("T7").Value = 0.09304561749543325
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.08683666666667
$ws.Range("H8").Value = 51.26051
$ws.Range("I8").Value = 0.1907597520636141
$ws.Range("J8").Value = 0.1907597520636141
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2909853333333334
$ws.Range("N8").Value = 0.8729560000000001
$ws.Range("O8").Value = 0.03275300561492853
$ws.Range("P8").Value = 0.03275300561492853
$ws.Range("Q8").Value = 4.972018863062223
$ws.Range("R8").Value = 44.74816976756001
$ws.Range("S8").Value = 0.006247955230441926
$ws.Range("T8").Value = 0.006247955230441926
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.08683666666667
$ws.Range("H9").Value = 51.26051
$ws.Range("I9").Value = 0.1907597520636141
$ws.Range("J9").Value = 0.1907597520636141
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.4814053333333333
$ws.Range("N9").Value = 1.444216
$ws.Range("O9").Value = 0.0541864822020464
$ws.Range("P9").Value = 0.05418648220204641
$ws.Range("Q9").Value = 8.225694301128888
$ws.Range("R9").Value = 74.03124871016
$ws.Range("S9").Value = 0.01033659991006181
$ws.Range("T9").Value = 0.01033659991006181
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 18.33915266666667
$ws.Range("H10").Value = 55.017458
$ws.Range("I10").Value = 0.2047407770084672
$ws.Range("J10").Value = 0.2047407770084672
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.778439
$ws.Range("N10").Value = 11.335317
$ws.Range("O10").Value = 0.4252971528324392
$ws.Range("P10").Value = 0.4252971528324392
$ws.Range("Q10").Value = 69.29336966268734
$ws.Range("R10").Value = 623.6403269641861
$ws.Range("S10").Value = 0.08707566953040242
$ws.Range("T10").Value = 0.08707566953040242
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 18.33915266666667
$ws.Range("H11").Value = 55.017458
$ws.Range("I11").Value = 0.2047407770084672
$ws.Range("J11").Value = 0.2047407770084672
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("N11").Value = 13.00021
$ws.Range("O11").Value = 0.4877633593505858
$ws.Range("P11").Value = 0.4877633593505858
$ws.Range("Q11").Value = 79.47094529624222
$ws.Range("R11").Value = 715.23850766618
$ws.Range("S11").Value = 0.09986504918969911
$ws.Range("T11").Value = 0.09986504918969913
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 18.33915266666667
$ws.Range("H12").Value = 55.017458
$ws.Range("I12").Value = 0.2047407770084672
$ws.Range("J12").Value = 0.2047407770084672
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.2909853333333334
$ws.Range("N12").Value = 0.8729560000000001
$ws.Range("O12").Value = 0.03275300561492853
$ws.Range("P12").Value = 0.03275300561492853
$ws.Range("Q12").Value = 5.33642445176089
$ws.Range("R12").Value = 48.027820065848
$ws.Range("S12").Value = 0.006705875818963155
$ws.Range("T12").Value = 0.006705875818963155
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 18.33915266666667
$ws.Range("H13").Value = 55.017458
$ws.Range("I13").Value = 0.2047407770084672
$ws.Range("J13").Value = 0.2047407770084672
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.4814053333333333
$ws.Range("N13").Value = 1.444216
$ws.Range("O13").Value = 0.0541864822020464
$ws.Range("P13").Value = 0.05418648220204641
$ws.Range("Q13").Value = 8.828565902547554
$ws.Range("R13").Value = 79.457093122928
$ws.Range("S13").Value = 0.01109418246940246
$ws.Range("T13").Value = 0.01109418246940246
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.438381
$ws.Range("H14").Value = 10.315143
$ws.Range("I14").Value = 0.03838654982521095
$ws.Range("J14").Value = 0.03838654982521095
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.778439
$ws.Range("N14").Value = 11.335317
$ws.Range("O14").Value = 0.4252971528324392
$ws.Range("P14").Value = 0.4252971528324392
$ws.Range("Q14").Value = 12.991712867259
$ws.Range("R14").Value = 116.925415805331
$ws.Range("S14").Value = 0.01632569034772279
$ws.Range("T14").Value = 0.01632569034772279
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.438381
$ws.Range("H15").Value = 10.315143
$ws.Range("I15").Value = 0.03838654982521095
$ws.Range("J15").Value = 0.03838654982521095
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("N15").Value = 13.00021
$ws.Range("O15").Value = 0.4877633593505858
$ws.Range("P15").Value = 0.4877633593505858
$ws.Range("Q15").Value = 14.89989168667
$ws.Range("R15").Value = 134.09902518003
$ws.Range("S15").Value = 0.01872355249662354
$ws.Range("T15").Value = 0.01872355249662354
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.438381
$ws.Range("H16").Value = 10.315143
$ws.Range("I16").Value = 0.03838654982521095
$ws.Range("J16").Value = 0.03838654982521095
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.2909853333333334
$ws.Range("N16").Value = 0.8729560000000001
$ws.Range("O16").Value = 0.03275300561492853
$ws.Range("P16").Value = 0.03275300561492853
$ws.Range("Q16").Value = 1.000518441412
$ws.Range("R16").Value = 9.004665972708002
$ws.Range("S16").Value = 0.001257274881962868
$ws.Range("T16").Value = 0.001257274881962868
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.438381
$ws.Range("H17").Value = 10.315143
$ws.Range("I17").Value = 0.03838654982521095
$ws.Range("J17").Value = 0.03838654982521095
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.4814053333333333
$ws.Range("N17").Value = 1.444216
$ws.Range("O17").Value = 0.0541864822020464
$ws.Range("P17").Value = 0.05418648220204641
$ws.Range("Q17").Value = 1.655254951432
$ws.Range("R17").Value = 14.897294562888
$ws.Range("S17").Value = 0.002080032098901761
$ws.Range("T17").Value = 0.002080032098901761
